$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ----------------------------------------------------------------------
# Finish metrics for the iceberg formats (rows 32-35).
# Columns: H=File Count, I=Record Count, J=Total Size (B), K=Avg File Size (B),
#          L=Table Name (already filled), M=Notes, N=Full Scan Time in Queue (ms),
#          O=Full Scan Run Time (sec), P=Data scanned (MB), Q=Job Run Time
# ----------------------------------------------------------------------

# Bring over the number-formatting (Comma style, time format, wrap-text Notes
# style, etc.) from an already fully populated row (29) so the new cells pick
# up the same styles already present in the workbook.
$ws.Range("H29:Q29").Copy()
$ws.Range("H32:Q32").PasteSpecial(-4122)

$ws.Range("H29:Q29").Copy()
$ws.Range("H33:Q33").PasteSpecial(-4122)
$ws.Range("I29").Copy()
$ws.Range("N33").PasteSpecial(-4122)

$ws.Range("H29:Q29").Copy()
$ws.Range("H34:Q34").PasteSpecial(-4122)

$ws.Range("H29:Q29").Copy()
$ws.Range("H35:Q35").PasteSpecial(-4122)
$ws.Range("I29:J29").Copy()
$ws.Range("N35:O35").PasteSpecial(-4122)

# Row heights grow because the Notes column text now wraps across more lines.
$ws.Rows.Item(32).RowHeight = 85.5
$ws.Rows.Item(33).RowHeight = 85.5
$ws.Rows.Item(34).RowHeight = 85.5
$ws.Rows.Item(35).RowHeight = 85.5

# --- Row 32: pricingtrxiceberggzip ---
$ws.Range("H32").Value = 12182
$ws.Range("I32").Value = 486906
$ws.Range("J32").Value = 73453103
$ws.Range("K32").Value = 6028
$ws.Range("M32").Value = "Fastest scan and smallest data size scanned. Puts all files in one S3 prefix s3://pricingtrxprocessed/iceberggzip/pricingtrx/pricingtrxiceberggzip/data/. Has a metadata prefix with 3 files. Files have .parquet suffix."
$ws.Range("N32").Value = 119
$ws.Range("O32").Value = 3.35
$ws.Range("P32").Value = 4.97
$ws.Range("Q32").Value = 0.018935185185185183

# --- Row 33: pricingtrxiceberglzo ---
$ws.Range("H33").Value = 12182
$ws.Range("I33").Value = 486906
$ws.Range("J33").Value = 73452937
$ws.Range("K33").Value = 6028
$ws.Range("M33").Value = "Fastest scan and smallest data size scanned. Puts all files in one S3 prefix s3://pricingtrxprocessed/iceberggzip/pricingtrx/pricingtrxiceberglzo/data/. Has a metadata prefix with 3 files. Files have .parquet suffix."
$ws.Range("N33").Value = 111
$ws.Range("O33").Value = 9.336
$ws.Range("P33").Value = 4.97
$ws.Range("Q33").Value = 0.018472222222222223

# --- Row 34: pricingtrxicebergsnappy ---
$ws.Range("H34").Value = 12182
$ws.Range("I34").Value = 486906
$ws.Range("J34").Value = 73454861
$ws.Range("K34").Value = 6028
$ws.Range("M34").Value = "Fastest scan and smallest data size scanned. Puts all files in one S3 prefix s3://pricingtrxprocessed/iceberggzip/pricingtrx/pricingtrxicebergsnappy/data/. Has a metadata prefix with 3 files. Files have .parquet suffix."
$ws.Range("N34").Value = 106
$ws.Range("O34").Value = 6.147
$ws.Range("P34").Value = 4.97
$ws.Range("Q34").Value = 0.01834490740740741

# --- Row 35: pricingtrxiceberguncompressed ---
$ws.Range("H35").Value = 12182
$ws.Range("I35").Value = 486906
$ws.Range("J35").Value = 73458508
$ws.Range("K35").Value = 6028
$ws.Range("M35").Value = "Fastest scan and smallest data size scanned. Puts all files in one S3 prefix s3://pricingtrxprocessed/iceberggzip/pricingtrx/pricingtrxiceberguncompressed/data/. Has a metadata prefix with 3 files. Files have .parquet suffix."
$ws.Range("N35").Value = 114
$ws.Range("O35").Value = 5.622
$ws.Range("P35").Value = 4.97
$ws.Range("Q35").Value = 0.017465277777777777

# Update the saved view/selection to match where the author ended up working.
$ws.Application.ActiveWindow.ScrollRow = 14
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("P35").Select()
